$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 288.94
$ws1.Range("M20").Value = 658.73
$ws1.Range("M35").Value = "4 de 33"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 288.94
$ws2.Range("F20").Value = 987.73
$ws2.Range("F35").Value = 9102.120000000001

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 7571.27
$ws3.Range("E12").Value = 14862.4853751766
$ws3.Range("F12").Value = 0.3374945421923323

$ws3.Range("D15").Value = 9291.310000000001
$ws3.Range("E15").Value = 29451.70881339593
$ws3.Range("F15").Value = 0.2398189476341839
